$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.794.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.46%  "

$ws.Range("E4").Value = "  -0.85%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.682"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.60%  "

$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.345"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0735"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0970"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.148.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.714"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.895.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "34.818.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0818"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "242.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.95%  "

$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("E25").Value = "  +4.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("E30").Value = "  -6.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0576"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.826"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.75%  "

$ws.Range("E39").Value = "  -23.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "97.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.282.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0797"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.94%  "

$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("E49").Value = "  -1.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.22%  "

